$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; existing rows 8.. shift down to 9..
$ws.Rows.Item(8).Insert()

# New row 8 mirrors the boilerplate columns of the (now shifted) row 9,
# with updated Fecha, Volumen, Precio mínimo/máximo/promedio and Precio $/Kg.
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = (Get-Date -Year 2023 -Month 10 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107002
$ws.Range("J8").Value = "Chirimoya"
$ws.Range("K8").Value = "Cultivar IV Región"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 23000
$ws.Range("P8").Value = 23000
$ws.Range("Q8").Value = "`$/bandeja 10 kilos"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 2300
$ws.Range("T8").Value = 10
